$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4645.7856
$ws.Range("I76").Value = 4442.375
$ws.Range("J76").Value = 4917
$ws.Range("K76").Value = 4442.375
$ws.Range("L76").Value = 4917
$ws.Range("M76").Value = -4127.375
$ws.Range("N76").Value = -5547
$ws.Range("H79").Value = 4645.7856
$ws.Range("I79").Value = 4442.375
$ws.Range("J79").Value = 4917
$ws.Range("K79").Value = 4442.375
$ws.Range("L79").Value = 4917
$ws.Range("M79").Value = -3350.375
$ws.Range("N79").Value = -7101
$ws.Range("H101").Value = 493.7143
$ws.Range("I101").Value = 495.57144
$ws.Range("J101").Value = 491.85715
$ws.Range("K101").Value = 1486.71432
$ws.Range("L101").Value = 1475.57145
$ws.Range("M101").Value = 135.28568
$ws.Range("N101").Value = -4719.571449999999
$ws.Range("H132").Value = 3382737.5
$ws.Range("I132").Value = 3792251
$ws.Range("J132").Value = 4250.625
$ws.Range("K132").Value = 11376753
$ws.Range("L132").Value = 12751.875
$ws.Range("M132").Value = -11374223
$ws.Range("N132").Value = -17811.875
$ws.Range("H137").Value = 1487.6389
$ws.Range("I137").Value = 1207.069
$ws.Range("K137").Value = 3621.207
$ws.Range("M137").Value = -1071.207
$ws.Range("H138").Value = 2621.5088
$ws.Range("I138").Value = 1326.7188
$ws.Range("J138").Value = 4278.84
$ws.Range("K138").Value = 3980.1564
$ws.Range("L138").Value = 12836.52
$ws.Range("M138").Value = 1159.8436
$ws.Range("N138").Value = -23116.52

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 28032
$ws.Range("I32").Value = 8512.93
$ws.Range("K32").Value = 8512.93
$ws.Range("M32").Value = -8225.93
$ws.Range("H63").Value = 2200
$ws.Range("J63").Value = 2200
$ws.Range("L63").Value = 2200
$ws.Range("N63").Value = -3572
$ws.Range("H66").Value = 2200
$ws.Range("J66").Value = 2200
$ws.Range("L66").Value = 11000
$ws.Range("N66").Value = -17864
$ws.Range("H102").Value = 52577.5
$ws.Range("I102").Value = 92845
$ws.Range("J102").Value = 3361.6667
$ws.Range("K102").Value = 92845
$ws.Range("L102").Value = 3361.6667
$ws.Range("M102").Value = -91223
$ws.Range("N102").Value = -6605.6667
$ws.Range("H122").Value = 2815.5293
$ws.Range("I122").Value = 2261.8462
$ws.Range("K122").Value = 6785.5386
$ws.Range("M122").Value = -4335.5386

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 287385.72
$ws.Range("I105").Value = 201937.8
$ws.Range("J105").Value = 501005.5
$ws.Range("K105").Value = 201937.8
$ws.Range("L105").Value = 501005.5
$ws.Range("M105").Value = -200190.8
$ws.Range("N105").Value = -504499.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 382.55554
$ws.Range("I19").Value = 382.55554
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 382.55554
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -212.55554
$ws.Range("N19").Value = $null
$ws.Range("H24").Value = 382.55554
$ws.Range("I24").Value = 382.55554
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 382.55554
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -212.55554
$ws.Range("N24").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 1800
$ws.Range("J9").Value = 1800
$ws.Range("L9").Value = 5400
$ws.Range("N9").Value = -5848
$ws.Range("H131").Value = 817217.3
$ws.Range("I131").Value = 680
$ws.Range("J131").Value = 938185.75
$ws.Range("K131").Value = 2040
$ws.Range("L131").Value = 2814557.25
$ws.Range("M131").Value = 3000
$ws.Range("N131").Value = -2824637.25
$ws.Range("H133").Value = 1113.25
$ws.Range("I133").Value = 1113.25
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 3339.75
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = 1720.25
$ws.Range("N133").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 62920.344
$ws.Range("I70").Value = 92225.914
$ws.Range("K70").Value = 92225.914
$ws.Range("M70").Value = -91955.914
$ws.Range("H73").Value = 62920.344
$ws.Range("I73").Value = 92225.914
$ws.Range("K73").Value = 92225.914
$ws.Range("M73").Value = -91289.914
$ws.Range("H80").Value = 167034990
$ws.Range("I80").Value = 250551250
$ws.Range("J80").Value = 2499.5
$ws.Range("K80").Value = 250551250
$ws.Range("L80").Value = 2499.5
$ws.Range("M80").Value = -250550252
$ws.Range("N80").Value = -4495.5
$ws.Range("H83").Value = 167034990
$ws.Range("I83").Value = 250551250
$ws.Range("J83").Value = 2499.5
$ws.Range("K83").Value = 1252756250
$ws.Range("L83").Value = 12497.5
$ws.Range("M83").Value = -1252751258
$ws.Range("N83").Value = -22481.5
$ws.Range("H126").Value = 2626.4211
$ws.Range("I126").Value = 2630.8
$ws.Range("K126").Value = 7892.400000000001
$ws.Range("M126").Value = -5422.400000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2154.6
$ws.Range("I61").Value = 2098.8
$ws.Range("J61").Value = 2210.4
$ws.Range("K61").Value = 2098.8
$ws.Range("L61").Value = 2210.4
$ws.Range("M61").Value = -1896.8
$ws.Range("N61").Value = -2614.4
$ws.Range("H68").Value = 3087.9714
$ws.Range("J68").Value = 6175.222
$ws.Range("L68").Value = 6175.222
$ws.Range("N68").Value = -7673.222
$ws.Range("H71").Value = 3087.9714
$ws.Range("J71").Value = 6175.222
$ws.Range("L71").Value = 30876.11
$ws.Range("N71").Value = -38364.11
$ws.Range("H113").Value = 2154.6
$ws.Range("I113").Value = 2098.8
$ws.Range("J113").Value = 2210.4
$ws.Range("K113").Value = 2098.8
$ws.Range("L113").Value = 2210.4
$ws.Range("M113").Value = 71.19999999999982
$ws.Range("N113").Value = -6550.4
$ws.Range("H132").Value = 4998.5
$ws.Range("I132").Value = 5798.4
$ws.Range("J132").Value = 999
$ws.Range("K132").Value = 17395.2
$ws.Range("L132").Value = 2997
$ws.Range("M132").Value = -14865.2
$ws.Range("N132").Value = -8057

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").Value = $null
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").Value = $null
$ws.Range("H114").Value = 32166.666
$ws.Range("J114").Value = 32166.666
$ws.Range("L114").Value = 32166.666
$ws.Range("N114").Value = -40844.666
